$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update serial_key (B) and username_pos (E) for each test-case row (TC_01..TC_46)
$ws.Range("B2").Value = "307260624WoJ"
$ws.Range("E2").Value = "usertwo_p1"
$ws.Range("B3").Value = "307260624WoJ"
$ws.Range("E3").Value = "usertwo_p1"
$ws.Range("B4").Value = "307260624WoJ"
$ws.Range("E4").Value = "usertwo_p1"
$ws.Range("B5").Value = "307260624WoJ"
$ws.Range("E5").Value = "usertwo_p1"
$ws.Range("B6").Value = "307260624WoJ"
$ws.Range("E6").Value = "usertwo_p1"
$ws.Range("B7").Value = "307260624WoJ"
$ws.Range("E7").Value = "usertwo_p1"
$ws.Range("B8").Value = "307260624WoJ"
$ws.Range("E8").Value = "usertwo_p1"
$ws.Range("B9").Value = "307260624WoJ"
$ws.Range("E9").Value = "usertwo_p1"
$ws.Range("B10").Value = "307260624WoJ"
$ws.Range("E10").Value = "usertwo_p1"
$ws.Range("B11").Value = "307260624WoJ"
$ws.Range("E11").Value = "usertwo_p1"
$ws.Range("B12").Value = "307260624WoJ"
$ws.Range("E12").Value = "usertwo_p1"
$ws.Range("B13").Value = "307260624WoJ"
$ws.Range("E13").Value = "usertwo_p1"
$ws.Range("B14").Value = "307260624WoJ"
$ws.Range("E14").Value = "usertwo_p1"
$ws.Range("B15").Value = "307260624WoJ"
$ws.Range("E15").Value = "usertwo_p1"
$ws.Range("B16").Value = "307260624WoJ"
$ws.Range("E16").Value = "usertwo_p1"
$ws.Range("B17").Value = "307260624WoJ"
$ws.Range("E17").Value = "usertwo_p1"
$ws.Range("B18").Value = "307260624WoJ"
$ws.Range("E18").Value = "usertwo_p1"
$ws.Range("B19").Value = "307260624WoJ"
$ws.Range("E19").Value = "usertwo_p1"
$ws.Range("B20").Value = "307260624WoJ"
$ws.Range("E20").Value = "usertwo_p1"
$ws.Range("B21").Value = "307260624WoJ"
$ws.Range("E21").Value = "usertwo_p1"
$ws.Range("B22").Value = "307260624WoJ"
$ws.Range("E22").Value = "usertwo_p1"
$ws.Range("B23").Value = "307260624WoJ"
$ws.Range("E23").Value = "usertwo_p1"
$ws.Range("B24").Value = "307260624WoJ"
$ws.Range("E24").Value = "usertwo_p1"
$ws.Range("B25").Value = "307260624WoJ"
$ws.Range("E25").Value = "usertwo_p1"
$ws.Range("B26").Value = "307260624ut0"
$ws.Range("E26").Value = "usertwo_p2"
$ws.Range("B27").Value = "307260624WoJ"
$ws.Range("E27").Value = "userone_p1"
$ws.Range("B28").Value = "307260624PTe"
$ws.Range("E28").Value = "userone_p4"
$ws.Range("B29").Value = "307260624uOm"
$ws.Range("E29").Value = "userone_p5"
$ws.Range("B30").Value = "307260624PTe"
$ws.Range("E30").Value = "userone_p4"
$ws.Range("B31").Value = "307260624uOm"
$ws.Range("E31").Value = "userone_p5"
$ws.Range("B32").Value = "307260624PTe"
$ws.Range("E32").Value = "userone_p4"
$ws.Range("B33").Value = "307260624PTe"
$ws.Range("E33").Value = "userone_p4"
$ws.Range("B34").Value = "307260624ut0"
$ws.Range("E34").Value = "userone_p2"
$ws.Range("B35").Value = "307260624WoJ"
$ws.Range("E35").Value = "usertwo_p1"
$ws.Range("B36").Value = "307260624SF2"
$ws.Range("E36").Value = "userone_p6"
$ws.Range("B37").Value = "307260624cuS"
$ws.Range("E37").Value = "userone_p7"
$ws.Range("B38").Value = "307260624SF2"
$ws.Range("E38").Value = "userone_p6"
$ws.Range("B39").Value = "307260624cuS"
$ws.Range("E39").Value = "userone_p7"
$ws.Range("B40").Value = "307260624SF2"
$ws.Range("E40").Value = "userone_p6"
$ws.Range("B41").Value = "307260624SF2"
$ws.Range("E41").Value = "userone_p6"
$ws.Range("B42").Value = "307260624WoJ"
$ws.Range("E42").Value = "usertwo_p1"
$ws.Range("B43").Value = "307260624WoJ"
$ws.Range("E43").Value = "usertwo_p1"
$ws.Range("B44").Value = "307260624WoJ"
$ws.Range("E44").Value = "usertwo_p1"
$ws.Range("B45").Value = "307260624WoJ"
$ws.Range("E45").Value = "usertwo_p1"
$ws.Range("B46").Value = "3072606245nJ"
$ws.Range("E46").Value = "userone_p3"
$ws.Range("B47").Value = "3072606245nJ"
$ws.Range("E47").Value = "userone_p3"

# TC_03 discount_value correction: Percentage_Item test case should use a percentage value
$ws.Range("Q4").Value = "Percentage : 10"

# Restore the view/selection state as saved by the author
$ws.Range("G13").Select()
